# Generate Report for Handback
# Refresh the timestamp columns that are re-stamped whenever the handback
# status report is (re)generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# "Latest HO Xliff Generate Date" for 91a92768-24d6-4634-8078-598cd5677dc5.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-31 17:21:02"

# --- zh-cn sheet -------------------------------------------------------
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the
# 91a92768-24d6-4634-8078-598cd5677dc5 zh-cn xliff
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-31 17:20:55"
$wsZhCn.Range("K2").Value = "2016-08-31 17:21:35"

# --- de-de sheet -------------------------------------------------------
# "Correspond Handback DateTime" for the 91a92768-24d6-4634-8078-598cd5677dc5
# de-de xliff
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-31 17:21:42"
